# Update the "dSF" column (F) values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -6
    4  = 2
    7  = -5
    10 = -2
    14 = 5
    17 = -2
    18 = 2
    21 = 2
    22 = 4
    25 = 3
    29 = -3
    30 = -1
    31 = 0
    32 = -2
    38 = 2
    39 = 1
    43 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
